$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

# Swap worker identity between row 16 and row 17, and update the "Periodo Mora"
# value from 2506 to 2507 for both rows.
$ws.Range("C16").Value = "1052072811"
$ws.Range("D16").Value = "SANDY LORENA BARRIOS YEPES"
$ws.Range("E16").Value = "2507"

$ws.Range("C17").Value = "1143334824"
$ws.Range("D17").Value = "YOMARIS PEROZA BERRIO"
$ws.Range("E17").Value = "2507"
